$wb = $excel.ActiveWorkbook

# Update data values on the "RAF-generation" sheet
$wsGen = $wb.Worksheets.Item("RAF-generation")
$wsGen.Range("B2").Value = 0.8    # hard coal
$wsGen.Range("B10").Value = 1     # geothermal
$wsGen.Range("B11").Value = 1     # petroleum
$wsGen.Range("B14").Value = 1     # lignite

# Move the active tab / selection from "About" to "RAF-generation"
$wsGen.Activate()
$wsGen.Range("B12").Select()

Write-Host "done"
